$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new SMS/transaction entry ("password internet") was logged after the most
# recent "Others" row (row 29) on 2024-09-03 20:05:31. Insert a fresh row so
# every later entry (including the "hdfc"/P-Q block and the "Broadband"
# group header) shifts down by one, then populate the new row's
# September_Details / September_Date columns (R/S).
$ws.Rows.Item(29).Insert()

$ws.Range("R29").Value = "password internet"
$ws.Range("S29").Value = "2024-09-03 20:05:31"
